$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new child "Ratty Teetherson" (sibling of row 6's "Chisel Teetherson") at Nutcracker Academy ---
$ws.Range("A7").Value = "Elementary School"
$ws.Range("B7").Value = "Nutcracker Academy"
$ws.Range("C7").Value = 1009876548
$ws.Range("D7").Value = "Ratty"
$ws.Range("E7").Value = "Teetherson"
$ws.Range("F7").Value = 17
$ws.Range("G7").NumberFormat = "mm-dd-yy"
$ws.Range("G7").Value2 = 39872
$ws.Range("J7").Value = "Hazelton"
$ws.Range("K7").Value = "Ontario"
$ws.Range("L7").Value = "N4U2L1"
$ws.Range("M7").Value = "HPV,"
$ws.Range("N7").Value = "HPV-9,"
$ws.Range("O7").Value = "Mar 12, 2014 - DTaP-IPV-Hib, Mar 12, 2014 - rota-unspecified, May 14, 2014 - Pneu-C-13, Jul 19, 2014 - DTaP-IPV-Hib, Sep 21, 2014 - MMR, Nov 25, 2014 - Men-C-C, Apr 17, 2015 - Var, Sep 13, 2015 - DTaP-IPV-Hib, May 5, 2024 - Tdap-IPV,"
$ws.Range("P7").Value = "NUTCRACKER ACADEMY-1009876547"
$ws.Range("Q7").Value = "HPV (HPV-9)"
$ws.Range("R7").Value = "[2014 MAR 12: DTaP-IPV-Hib, rota-unspecified] [2014 MAY 14: Pneu-C-13] [2014 JUL 19: DTaP-IPV-Hib] [2014 SEP 21: MMR] [2014 NOV 25: Men-C-C] [2015 APR 17: Var] [2015 SEP 13: DTaP-IPV-Hib] [2024 MAY 05: Tdap-IPV]"

# --- Row 8: new child "Cheddarina Swiftpaws" at Tunnel Academy ---
$ws.Range("A8").Value = "Elementary School"
$ws.Range("B8").Value = "Tunnel Academy"
$ws.Range("C8").Value = 1009876550
$ws.Range("D8").Value = "Cheddarina"
$ws.Range("E8").Value = "Swiftpaws"
$ws.Range("F8").Value = 11
[void]$ws.Range("G7").Copy()
[void]$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value2 = 41896
$ws.Range("H8").Value = "44 Hayloft Road"
$ws.Range("J8").Value = "Burrowville"
$ws.Range("K8").Value = "Ontario"
$ws.Range("L8").Value = "H8Y6T6"
$ws.Range("M8").Value = "MMR"
$ws.Range("N8").Value = "MMR"
$ws.Range("P8").Value = "TUNNEL ACADEMY-1009876550"
$ws.Range("O8").Value = "Jan 10, 2015 - DTaP-IPV-Hib, Jan 29, 2015 - Pneu-C-13, Feb 18, 2015 - rota-unspecified, Mar 07, 2015 - DTaP-IPV-Hib, Mar 28, 2015 - MMR, Apr 15, 2015 - Men-C-C, May 02, 2015 - Var, May 27, 2015 - DTaP-IPV-Hib, Jun 16, 2015 - Pneu-C-13, Jul 09, 2015 - Influenza (IIV4), Aug 01, 2015 - Influenza (IIV4), Aug 29, 2015 - MMR, Sep 22, 2015 - Var, Oct 11, 2015 - DTaP-IPV-Hib, Nov 05, 2015 - Pneu-C-13, Dec 03, 2015 - Men-C-C, Jan 14, 2016 - MMR, Feb 06, 2016 - Influenza (IIV4), Mar 12, 2016 - Hep A, Apr 04, 2016 - Hep A booster, May 18, 2016 - Yellow Fever, Jun 07, 2016 - Rabies (pre-exposure), Jun 30, 2016 - Rabies (pre-exposure) dose 2, Jul 23, 2016 - Rabies (pre-exposure) dose 3, Aug 15, 2016 - Var, Sep 08, 2016 - DTaP-IPV-Hib, Oct 01, 2016 - Pneu-C-13, Oct 27, 2016 - Influenza (IIV4), Nov 19, 2016 - MMR, Dec 14, 2016 - Men-C-C, Jan 09, 2017 - Var, Feb 03, 2017 - DTaP-IPV-Hib, Mar 01, 2017 - Pneu-C-13, Mar 29, 2017 - MMR, Apr 18, 2017 - Influenza (IIV4), May 10, 2017 - COVID-19 (Pfizer Pediatric), Jun 02, 2017 - COVID-19 (Pfizer Pediatric) dose 2, Jun 28, 2017 - COVID-19 Booster, Jul 20, 2017 - Var, Aug 12, 2017 - Men-C-C, Sep 03, 2017 - Influenza (IIV4), Oct 25, 2017 - DTaP-IPV-Hib, Nov 16, 2017 - Pneu-C-13, Dec 08, 2017 - MMR, May 02, 2023 - Tdap, Jan 18, 2024 - Men-C-ACYW-135, May 01, 2024 - Tdap-IPV"
$ws.Range("Q8").Value = "Measles (MMR)"
$ws.Range("R8").Value = "[2015 JAN 10: DTaP-IPV-Hib] [2015 JAN 29: Pneu-C-13] [2015 FEB 18: rota-unspecified] [2015 MAR 07: DTaP-IPV-Hib] [2015 MAR 28: MMR] [2015 APR 15: Men-C-C] [2015 MAY 02: Var] [2015 MAY 27: DTaP-IPV-Hib] [2015 JUN 16: Pneu-C-13] [2015 JUL 09: Influenza (IIV4)] [2015 AUG 01: Influenza (IIV4)] [2015 AUG 29: MMR] [2015 SEP 22: Var] [2015 OCT 11: DTaP-IPV-Hib] [2015 NOV 05: Pneu-C-13] [2015 DEC 03: Men-C-C] [2016 JAN 14: MMR] [2016 FEB 06: Influenza (IIV4)] [2016 MAR 12: Hep A] [2016 APR 04: Hep A booster] [2016 MAY 18: Yellow Fever] [2016 JUN 07: Rabies (pre-exposure)] [2016 JUN 30: Rabies (pre-exposure) dose 2] [2016 JUL 23: Rabies (pre-exposure) dose 3] [2016 AUG 15: Var] [2016 SEP 08: DTaP-IPV-Hib] [2016 OCT 01: Pneu-C-13] [2016 OCT 27: Influenza (IIV4)] [2016 NOV 19: MMR] [2016 DEC 14: Men-C-C] [2017 JAN 09: Var] [2017 FEB 03: DTaP-IPV-Hib] [2017 MAR 01: Pneu-C-13] [2017 MAR 29: MMR] [2017 APR 18: Influenza (IIV4)] [2017 MAY 10: COVID-19 (Pfizer Pediatric)] [2017 JUN 02: COVID-19 (Pfizer Pediatric) dose 2] [2017 JUN 28: COVID-19 Booster] [2017 JUL 20: Var] [2017 AUG 12: Men-C-C] [2017 SEP 03: Influenza (IIV4)] [2017 OCT 25: DTaP-IPV-Hib] [2017 NOV 16: Pneu-C-13] [2017 DEC 08: MMR] [2023 MAY 02: Tdap] [2024 JAN 18: Men-C-ACYW-135] [2024 MAY 01: Tdap-IPV]"

# --- Column G width to fit the new date values (AutoFit, matches real-Excel behavior of widening a date column) ---
[void]$ws.Columns.Item(7).AutoFit()

# --- Leave the same selection state Excel ends up in after entering the last row ---
[void]$ws.Range("A8:R8").Select()
